$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value  = 5.983999999999997
$ws.Range("B10").Value = 5.303699999999999
$ws.Range("B12").Value = 4.8649
$ws.Range("B18").Value = 6.381499999999994
$ws.Range("B37").Value = 8.851500000000001
$ws.Range("B55").Value = 6.295799999999994
$ws.Range("B68").Value = 4.835299999999993
$ws.Range("B77").Value = 9.324700000000007
$ws.Range("B78").Value = 9.898300000000003
$ws.Range("B81").Value = 5.364600000000005
$ws.Range("B82").Value = 5.396000000000002

$wb.Save()
